$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 166668180
$ws.Range("I32").Value = 333334560
$ws.Range("J32").Value = 1794
$ws.Range("K32").Value = 333334560
$ws.Range("L32").Value = 1794
$ws.Range("M32").Value = -333334234
$ws.Range("N32").Value = -2446
$ws.Range("H94").Value = 9335.833000000001
$ws.Range("I94").Value = 3232.8572
$ws.Range("K94").Value = 3232.8572
$ws.Range("M94").Value = -2781.8572
$ws.Range("H116").Value = 3666.1836
$ws.Range("I116").Value = 3484.2144
$ws.Range("J116").Value = 3908.8096
$ws.Range("K116").Value = 3484.2144
$ws.Range("L116").Value = 3908.8096
$ws.Range("M116").Value = -42.21439999999984
$ws.Range("N116").Value = -10792.8096
$ws.Range("H125").Value = 1390911
$ws.Range("I125").Value = 990
$ws.Range("J125").Value = 1517267.5
$ws.Range("K125").Value = 8910
$ws.Range("L125").Value = 13655407.5
$ws.Range("M125").Value = -6450
$ws.Range("N125").Value = -13660327.5
$ws.Range("H128").Value = 33666.668
$ws.Range("J128").Value = 33666.668
$ws.Range("L128").Value = 33666.668
$ws.Range("N128").Value = -43626.668
$ws.Range("H129").Value = 706.3333
$ws.Range("I129").Value = 407.86365
$ws.Range("J129").Value = 2019.6
$ws.Range("K129").Value = 1223.59095
$ws.Range("L129").Value = 6058.799999999999
$ws.Range("M129").Value = 3776.40905
$ws.Range("N129").Value = -16058.8
$ws.Range("H132").Value = 2020.8195
$ws.Range("I132").Value = 1217.4237
$ws.Range("J132").Value = 5667
$ws.Range("K132").Value = 3652.2711
$ws.Range("L132").Value = 17001
$ws.Range("M132").Value = -1122.2711
$ws.Range("N132").Value = -22061
$ws.Range("H138").Value = 2533.3635
$ws.Range("I138").Value = 1935.7778
$ws.Range("J138").Value = 3250.4666
$ws.Range("K138").Value = 5807.3334
$ws.Range("L138").Value = 9751.399800000001
$ws.Range("N138").Value = -20031.3998
$ws.Range("M138").Value = -667.3334000000004

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3494860.8
$ws.Range("I32").Value = 5784.8623
$ws.Range("K32").Value = 5784.8623
$ws.Range("M32").Value = -5497.8623
$ws.Range("H36").Value = 166668460
$ws.Range("I36").Value = 202.5
$ws.Range("J36").Value = 500005000
$ws.Range("K36").Value = 202.5
$ws.Range("L36").Value = 500005000
$ws.Range("M36").Value = 143.5
$ws.Range("N36").Value = -500005692
$ws.Range("H45").Value = 4400.3687
$ws.Range("I45").Value = 3938.875
$ws.Range("J45").Value = 4736
$ws.Range("K45").Value = 3938.875
$ws.Range("L45").Value = 4736
$ws.Range("M45").Value = -3561.875
$ws.Range("N45").Value = -5490
$ws.Range("H94").Value = 41443.332
$ws.Range("J94").Value = 41443.332
$ws.Range("L94").Value = 41443.332
$ws.Range("N94").Value = -43245.332

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 30379.572
$ws.Range("J103").Value = 30379.572
$ws.Range("L103").Value = 30379.572
$ws.Range("N103").Value = -32723.572
$ws.Range("H107").Value = 5486.875
$ws.Range("I107").Value = 6139.5
$ws.Range("J107").Value = 3529
$ws.Range("K107").Value = 6139.5
$ws.Range("L107").Value = 3529
$ws.Range("M107").Value = -4219.5
$ws.Range("N107").Value = -7369

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 8233.333000000001
$ws.Range("J43").Value = 8233.333000000001
$ws.Range("L43").Value = 8233.333000000001
$ws.Range("N43").Value = -8601.333000000001
$ws.Range("H99").Value = 2274.75
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 2799.5
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 2799.5
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -5795.5
$ws.Range("H101").Value = 8233.333000000001
$ws.Range("J101").Value = 8233.333000000001
$ws.Range("L101").Value = 8233.333000000001
$ws.Range("N101").Value = -14723.333
$ws.Range("H126").Value = 2274.75
$ws.Range("I126").Value = 1750
$ws.Range("J126").Value = 2799.5
$ws.Range("K126").Value = 5250
$ws.Range("L126").Value = 8398.5
$ws.Range("M126").Value = -2780
$ws.Range("N126").Value = -13338.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 363.82352
$ws.Range("I12").Value = 58.555557
$ws.Range("J12").Value = 473.72
$ws.Range("K12").Value = 175.666671
$ws.Range("L12").Value = 1421.16
$ws.Range("M12").Value = -2.666671000000008
$ws.Range("N12").Value = -1767.16
$ws.Range("H70").Value = 6338.4
$ws.Range("I70").Value = 5584.3335
$ws.Range("J70").Value = 7034.4614
$ws.Range("K70").Value = 16753.0005
$ws.Range("L70").Value = 21103.3842
$ws.Range("M70").Value = -16438.0005
$ws.Range("N70").Value = -21733.3842
$ws.Range("H73").Value = 6338.4
$ws.Range("I73").Value = 5584.3335
$ws.Range("J73").Value = 7034.4614
$ws.Range("K73").Value = 16753.0005
$ws.Range("L73").Value = 21103.3842
$ws.Range("M73").Value = -15661.0005
$ws.Range("N73").Value = -23287.3842
$ws.Range("H107").Value = 656.5769
$ws.Range("I107").Value = 330.55
$ws.Range("J107").Value = 1743.3334
$ws.Range("K107").Value = 991.6500000000001
$ws.Range("L107").Value = 5230.0002
$ws.Range("M107").Value = 928.3499999999999
$ws.Range("N107").Value = -9070.0002
$ws.Range("H112").Value = 29799374
$ws.Range("J112").Value = 33375186
$ws.Range("L112").Value = 100125558
$ws.Range("N112").Value = -100127774
$ws.Range("H122").Value = 50894.227
$ws.Range("I122").Value = 352
$ws.Range("J122").Value = 57724.258
$ws.Range("K122").Value = 3168
$ws.Range("L122").Value = 519518.322
$ws.Range("M122").Value = -718
$ws.Range("N122").Value = -524418.322
$ws.Range("H129").Value = 1333.871
$ws.Range("I129").Value = 735.38464
$ws.Range("J129").Value = 1766.1111
$ws.Range("K129").Value = 2206.15392
$ws.Range("L129").Value = 5298.3333
$ws.Range("M129").Value = 2793.84608
$ws.Range("N129").Value = -15298.3333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 12525.2
$ws.Range("J134").Value = 12525.2
$ws.Range("L134").Value = 37575.60000000001
$ws.Range("N134").Value = -42645.60000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 75000
$ws.Range("J103").Value = 75000
$ws.Range("L103").Value = 75000
$ws.Range("N103").Value = -77344

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 370944.88
$ws.Range("I100").Value = 363.93332
$ws.Range("J100").Value = 834171.0600000001
$ws.Range("K100").Value = 727.86664
$ws.Range("L100").Value = 1668342.12
$ws.Range("M100").Value = -186.86664
$ws.Range("N100").Value = -1669424.12
$ws.Range("H101").Value = 4600.6665
$ws.Range("J101").Value = 4600.6665
$ws.Range("L101").Value = 4600.6665
$ws.Range("N101").Value = -11090.6665
$ws.Range("H136").Value = 1310.2428
$ws.Range("I136").Value = 1234.4329
$ws.Range("K136").Value = 3703.2987
$ws.Range("M136").Value = -1153.2987
